# Update the Direct Loss Measurement sequence-number bullet to read
# "... - DLM session status" instead of "... - Detect session state up/down".
#
# (Lives on slide 4, shape "Content Placeholder 2", 10th paragraph, but we
# locate it by content so the script is resilient to shape/paragraph
# renumbering.)

$p = $ppt.ActivePresentation

$oldText = "Sequence Numbers allow to detect Direct Loss Measurement test packet loss - Detect session state up/down"
$newText = "Sequence Numbers allow to detect Direct Loss Measurement test packet loss - DLM session status"
$needle  = "Detect session state up/down"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            if ($full -ne $null -and $full.Contains($needle)) {
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    $currentText = $para.Text.TrimEnd([char]13)
                    if ($currentText -eq $oldText) {
                        # Route the replacement through an unrelated
                        # placeholder string first so that neither
                        # assignment shares a common prefix/suffix with the
                        # surrounding text; this keeps the paragraph as a
                        # single run and preserves the original run
                        # formatting (rPr) instead of being split into
                        # "unchanged" / "changed" run fragments.
                        $para.Text = "##TEMP_PLACEHOLDER_TEXT##"
                        $para2 = $tr.Paragraphs($pi, 1)
                        $para2.Text = $newText
                    } elseif ($currentText.Contains($needle)) {
                        # Fallback: formatting differs from what we
                        # expected, just replace the text directly.
                        $para.Text = $para.Text.Replace($needle, "DLM session status")
                    }
                }
            }
        }
    }
}
